$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.064803118933699
$ws.Range("D2").Value = 1.063342780139254
$ws.Range("E2").Value = 1.068114835142486
$ws.Range("F2").Value = 1.074854466764149
$ws.Range("I2").Value = 1.028648454665347
$ws.Range("J2").Value = 1.069761488006837
$ws.Range("K2").Value = 1.066061677507009
$ws.Range("L2").Value = 1.070820883502488
$ws.Range("M2").Value = 1.077542578388205
$ws.Range("N2").Value = 1.071280672023164
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.067841440911353
$ws.Range("D3").Value = 1.06612893956993
$ws.Range("E3").Value = 1.070801267243663
$ws.Range("F3").Value = 1.077565724272627
$ws.Range("I3").Value = 1.028732841268858
$ws.Range("J3").Value = 1.072445289606057
$ws.Range("K3").Value = 1.068657890154306
$ws.Range("L3").Value = 1.073318588937215
$ws.Range("M3").Value = 1.080066402219196
$ws.Range("N3").Value = 1.073968284928492
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.069794240581435
$ws.Range("D4").Value = 1.067919331974356
$ws.Range("E4").Value = 1.072527093475551
$ws.Range("F4").Value = 1.079306831570139
$ws.Range("I4").Value = 1.028783977858426
$ws.Range("J4").Value = 1.074168962578275
$ws.Range("K4").Value = 1.070325168478559
$ws.Range("L4").Value = 1.074922064028128
$ws.Range("M4").Value = 1.081685995128424
$ws.Range("N4").Value = 1.075694405714038
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.070612127963329
$ws.Range("D5").Value = 1.068669116758371
$ws.Range("E5").Value = 1.073249722875813
$ws.Range("F5").Value = 1.080035697784476
$ws.Range("I5").Value = 1.028804650050469
$ws.Range("J5").Value = 1.074890580635758
$ws.Range("K5").Value = 1.071023143217264
$ws.Range("L5").Value = 1.075593197972291
$ws.Range("M5").Value = 1.082363717647209
$ws.Range("N5").Value = 1.076417048551934
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.070749277434424
$ws.Range("D6").Value = 1.068794841436469
$ws.Range("E6").Value = 1.073370887264649
$ws.Range("F6").Value = 1.080157898376158
$ws.Range("I6").Value = 1.028808072713358
$ws.Range("J6").Value = 1.075011569047491
$ws.Range("K6").Value = 1.071140165555577
$ws.Range("L6").Value = 1.075705712487539
$ws.Range("M6").Value = 1.082477327425824
$ws.Range("N6").Value = 1.076538208781089
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.069805181175169
$ws.Range("D7").Value = 1.067929361907283
$ws.Range("E7").Value = 1.072536760607298
$ws.Range("F7").Value = 1.079316582774577
$ws.Range("I7").Value = 1.028784257319507
$ws.Range("J7").Value = 1.07417861661088
$ws.Range("K7").Value = 1.070334506335707
$ws.Range("L7").Value = 1.074931043311548
$ws.Range("M7").Value = 1.081695063184862
$ws.Range("N7").Value = 1.075704073456477
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.065832721540007
$ws.Range("D8").Value = 1.064287001113999
$ws.Range("E8").Value = 1.069025358312096
$ws.Range("F8").Value = 1.075773541859688
$ws.Range("I8").Value = 1.028677694121013
$ws.Range("J8").Value = 1.070671216753163
$ws.Range("K8").Value = 1.066941745493738
$ws.Range("L8").Value = 1.071667671653904
$ws.Range("M8").Value = 1.078398355395186
$ws.Range("N8").Value = 1.072191692688657
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.058727539156332
$ws.Range("D9").Value = 1.057769683307074
$ws.Range("E9").Value = 1.06273867565982
$ws.Range("F9").Value = 1.069425112458935
$ws.Range("I9").Value = 1.028463156652967
$ws.Range("J9").Value = 1.064388044721041
$ws.Range("K9").Value = 1.060862843502671
$ws.Range("L9").Value = 1.065816426539951
$ws.Range("M9").Value = 1.072482374613625
$ws.Range("N9").Value = 1.06589959783156
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.053914142227469
$ws.Range("D10").Value = 1.053352894200224
$ws.Range("E10").Value = 1.058475761367307
$ws.Range("F10").Value = 1.065116999489558
$ws.Range("I10").Value = 1.028301843409431
$ws.Range("J10").Value = 1.060124965100927
$ws.Range("K10").Value = 1.056737617463323
$ws.Range("L10").Value = 1.061842932968151
$ws.Range("M10").Value = 1.068461693385813
$ws.Range("N10").Value = 1.061630464149403
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.051810375036396
$ws.Range("D11").Value = 1.051422099653236
$ws.Range("E11").Value = 1.056611671158382
$ws.Range("F11").Value = 1.063232379842366
$ws.Range("I11").Value = 1.028227586195438
$ws.Range("J11").Value = 1.058260180539964
$ws.Range("K11").Value = 1.054932956951989
$ws.Range("L11").Value = 1.060104008802551
$ws.Range("M11").Value = 1.066701379836958
$ws.Range("N11").Value = 1.05976303138047
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.051025889532612
$ws.Range("D12").Value = 1.050702060779343
$ws.Range("E12").Value = 1.055916424686981
$ws.Range("F12").Value = 1.062529363749193
$ws.Range("I12").Value = 1.028199335389026
$ws.Range("J12").Value = 1.057564580091119
$ws.Range("K12").Value = 1.054259757307709
$ws.Range("L12").Value = 1.059455235668309
$ws.Range("M12").Value = 1.06604451789037
$ws.Range("N12").Value = 1.059066443099201
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.051194304408946
$ws.Range("D13").Value = 1.050856642548641
$ws.Range("E13").Value = 1.056065687635945
$ws.Range("F13").Value = 1.062680299830267
$ws.Range("I13").Value = 1.028205425630045
$ws.Range("J13").Value = 1.057713923357612
$ws.Range("K13").Value = 1.054404292395569
$ws.Range("L13").Value = 1.059594530712313
$ws.Range("M13").Value = 1.066185554522896
$ws.Range("N13").Value = 1.059215998450262
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.051745592188712
$ws.Range("D14").Value = 1.051362639904131
$ws.Range("E14").Value = 1.056554260417879
$ws.Range("F14").Value = 1.063174329802368
$ws.Range("I14").Value = 1.028225264646178
$ws.Range("J14").Value = 1.058202742514026
$ws.Range("K14").Value = 1.054877369177122
$ws.Range("L14").Value = 1.060050439922406
$ws.Range("M14").Value = 1.066647145267998
$ws.Range("N14").Value = 1.059705511785946
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.052084850687913
$ws.Range("D15").Value = 1.05167401989537
$ws.Range("E15").Value = 1.056854906724887
$ws.Range("F15").Value = 1.063478319369055
$ws.Range("I15").Value = 1.02823739936414
$ws.Range("J15").Value = 1.058503527814859
$ws.Range("K15").Value = 1.055168464192695
$ws.Range("L15").Value = 1.060330958753088
$ws.Range("M15").Value = 1.066931145616353
$ws.Range("N15").Value = 1.060006724236408
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.054053340775595
$ws.Range("D16").Value = 1.053480640136242
$ws.Range("E16").Value = 1.058599082294152
$ws.Range("F16").Value = 1.06524166269054
$ws.Range("I16").Value = 1.02830667823147
$ws.Range("J16").Value = 1.060248318777045
$ws.Range("K16").Value = 1.056856990277479
$ws.Range("L16").Value = 1.061957944142286
$ws.Range("M16").Value = 1.068578103936241
$ws.Range("N16").Value = 1.061753993001891
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.055282810522547
$ws.Range("D17").Value = 1.054608912377673
$ws.Range("E17").Value = 1.059688207717183
$ws.Range("F17").Value = 1.066342554936718
$ws.Range("I17").Value = 1.028348950716307
$ws.Range("J17").Value = 1.061337661747505
$ws.Range("K17").Value = 1.057911157719747
$ws.Range("L17").Value = 1.062973520813787
$ws.Range("M17").Value = 1.069605952871002
$ws.Range("N17").Value = 1.062844882964318
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.055998062753615
$ws.Range("D18").Value = 1.055265256602982
$ws.Range("E18").Value = 1.060321725270748
$ws.Range("F18").Value = 1.066982843444614
$ws.Range("I18").Value = 1.028373182710125
$ws.Range("J18").Value = 1.061971246277926
$ws.Range("K18").Value = 1.058524266481536
$ws.Range("L18").Value = 1.063564123363702
$ws.Range("M18").Value = 1.070203621312534
$ws.Range("N18").Value = 1.063479367257445
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.05624163048534
$ws.Range("D19").Value = 1.055488758053856
$ws.Range("E19").Value = 1.060537444482286
$ws.Range("F19").Value = 1.067200855701783
$ws.Range("I19").Value = 1.028381373319228
$ws.Range("J19").Value = 1.062186978238866
$ws.Range("K19").Value = 1.058733023685605
$ws.Range("L19").Value = 1.063765206948216
$ws.Range("M19").Value = 1.070407098528353
$ws.Range("N19").Value = 1.063695405582514
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.055151095008557
$ws.Range("D20").Value = 1.054488042034773
$ws.Range("E20").Value = 1.059571536703624
$ws.Range("F20").Value = 1.066224631016253
$ws.Range("I20").Value = 1.028344459265166
$ws.Range("J20").Value = 1.061220973608263
$ws.Range("K20").Value = 1.057798239265061
$ws.Range("L20").Value = 1.062864742432393
$ws.Range("M20").Value = 1.069495867340812
$ws.Range("N20").Value = 1.062728029114534
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.051583336910694
$ws.Range("D21").Value = 1.05121371600412
$ws.Range("E21").Value = 1.056410467160771
$ws.Range("F21").Value = 1.063028933517012
$ws.Range("I21").Value = 1.028219441051515
$ws.Range("J21").Value = 1.058058879378084
$ws.Range("K21").Value = 1.054738139874437
$ws.Range("L21").Value = 1.059916265759108
$ws.Range("M21").Value = 1.066511302005811
$ws.Range("N21").Value = 1.059561444347851
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049322422938154
$ws.Range("D22").Value = 1.049138439120744
$ws.Range("E22").Value = 1.054406490579207
$ws.Range("F22").Value = 1.061002351323617
$ws.Range("I22").Value = 1.028136966922653
$ws.Range("J22").Value = 1.056053701440421
$ws.Range("K22").Value = 1.052797485476316
$ws.Range("L22").Value = 1.058045847638856
$ws.Range("M22").Value = 1.064617361226256
$ws.Range("N22").Value = 1.057553418827519
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.050522696389608
$ws.Range("D23").Value = 1.050240190654693
$ws.Range("E23").Value = 1.055470434407311
$ws.Range("F23").Value = 1.062078357879606
$ws.Range("I23").Value = 1.028181056996184
$ws.Range("J23").Value = 1.057118335792264
$ws.Range("K23").Value = 1.053827876156322
$ws.Range("L23").Value = 1.059038998255858
$ws.Range("M23").Value = 1.065623060574729
$ws.Range("N23").Value = 1.058619565082257
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.05521061737186
$ws.Range("D24").Value = 1.05454266356568
$ws.Range("E24").Value = 1.059624260716971
$ws.Range("F24").Value = 1.066277921449956
$ws.Range("I24").Value = 1.028346490071037
$ws.Range("J24").Value = 1.06127370554493
$ws.Range("K24").Value = 1.057849267715412
$ws.Range("L24").Value = 1.062913900152391
$ws.Range("M24").Value = 1.069545615990366
$ws.Range("N24").Value = 1.062780835936599
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.060577461179434
$ws.Range("D25").Value = 1.059466840484667
$ws.Range("E25").Value = 1.064376202729487
$ws.Range("F25").Value = 1.071079313321076
$ws.Range("I25").Value = 1.028521820747061
$ws.Range("J25").Value = 1.066025096535231
$ws.Range("K25").Value = 1.062446803120786
$ws.Range("L25").Value = 1.067341552620832
$ws.Range("M25").Value = 1.074024946704704
$ws.Range("N25").Value = 1.067538974447099
